$d = $word.ActiveDocument

$replacements = @(
    @("45×43=", "64×53="),
    @("76×72=", "62×47="),
    @("86×59=", "34×43="),
    @("48×73=", "36×87="),
    @("16×58=", "89×27="),
    @("37×95=", "25×94="),
    @("38×51=", "80×97="),
    @("18×84=", "75×79="),
    @("77×93=", "79×48="),
    @("49×73=", "40×94="),
    @("90×70=", "24×37="),
    @("73×15=", "33×49="),
    @("72×65=", "72×60="),
    @("45×91=", "88×95="),
    @("38×95=", "66×61="),
    @("80×67=", "68×59="),
    @("76×97=", "17×46="),
    @("52×56=", "87×39="),
    @("17×75=", "67×56="),
    @("28×68=", "57×31="),
    @("38×69=", "67×68="),
    @("65×46=", "22×77="),
    @("29×41=", "42×64="),
    @("53×23=", "69×79="),
    @("11×76=", "18×44=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
